# "tambah export di kelompok"
# The "DESA" label in the "kelompok" (group) export template is being
# replaced with a lower-case "desa" value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "desa"

# Reset the selection back to the default top-left cell (the stale
# A5:A7 selection from the source template no longer applies).
$ws.Range("A1").Select() | Out-Null
